$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    if ($null -eq $text) { return }
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Simple numeric/percentage updates: row -> D (price text), E (volume % text)
$updates = @(
    @{ Row = 2;  D = "87.311.59";   E = "  -3.48%  " },
    @{ Row = 3;  D = "3.055.59";    E = "  -4.28%  " },
    @{ Row = 4;  D = "0.999";       E = "  -0.22%  " },
    @{ Row = 5;  D = "208.96";      E = "  -2.87%  " },
    @{ Row = 6;  D = "621.00";      E = "  +0.10%  " },
    @{ Row = 7;  D = "0.360";       E = "  -10.21%  " },
    @{ Row = 8;  D = "0.769";       E = "  +11.09%  " },
    @{ Row = 9;  D = "0.999";       E = "  -0.03%  " },
    @{ Row = 10; D = "3.048.56";    E = "  -4.35%  " },
    @{ Row = 11; D = "0.576";       E = "  -0.47%  " },
    @{ Row = 12; D = $null;         E = "  -0.50%  " },
    @{ Row = 13; D = "0.0000233";   E = "  -10.23%  " },
    @{ Row = 14; D = "5.22";        E = "  -0.89%  " },
    @{ Row = 15; D = "87.077.18";   E = "  -3.51%  " },
    @{ Row = 16; D = $null;         E = "  -4.59%  " },
    @{ Row = 17; D = "31.12";       E = "  -6.22%  " },
    @{ Row = 18; D = "3.057.15";    E = "  -3.32%  " },
    @{ Row = 19; D = "3.34";        E = "  +1.87%  " },
    @{ Row = 20; D = "0.0000207";   E = "  -1.49%  " },
    @{ Row = 21; D = "12.94";       E = "  -3.99%  " },
    @{ Row = 22; D = "413.70";      E = "  -5.93%  " },
    @{ Row = 23; D = "8.14";        E = "  -5.75%  " },
    @{ Row = 24; D = "4.76";        E = "  -6.54%  " },
    @{ Row = 25; D = "5.37";        E = "  +3.83%  " },
    @{ Row = 26; D = "81.79";       E = "  +8.28%  " },
    @{ Row = 27; D = "11.16";       E = "  -4.53%  " },
    @{ Row = 28; D = "3.215.83";    E = "  -4.12%  " },
    @{ Row = 29; D = $null;         E = "  +0.22%  " },
    @{ Row = 30; D = "1.00";        E = "  +0.16%  " },
    @{ Row = 31; D = "0.149";       E = "  -12.54%  " },
    @{ Row = 32; D = "7.95";        E = "  -6.55%  " },
    @{ Row = 33; D = "492.48";      E = "  -8.54%  " },
    @{ Row = 34; D = "3.60";        E = "  -14.31%  " },
    @{ Row = 35; D = $null;         E = "  +8.70%  " },
    @{ Row = 36; D = "6.57";        E = "  -6.60%  " },
    @{ Row = 37; D = "1.77";        E = "  -5.31%  " },
    @{ Row = 38; D = "1.23";        E = "  -3.21%  " },
    @{ Row = 41; D = $null;         E = "  +0.14%  " },
    @{ Row = 43; D = "0.358";       E = "  -4.86%  " },
    @{ Row = 48; D = "0.0635";      E = "  +6.55%  " },
    @{ Row = 49; D = "157.42";      E = "  -9.06%  " },
    @{ Row = 50; D = "0.698";       E = "  -1.31%  " },
    @{ Row = 51; D = "1.16";        E = "  -6.80%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    Set-TextCell $r 4 $u.D
    Set-TextCell $r 5 $u.E
}

# Rows where the coin (B), link (C), price (D) and volume (E) got swapped/reordered
$rowSwaps = @(
    @{ Row = 39; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "21.73";  E = "  -1.74%  " },
    @{ Row = 40; B = "WhiteBITCoin";     C = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt";         D = "22.13";  E = "  -1.11%  " },
    @{ Row = 44; B = "Stacks";           C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx";                D = "1.80";   E = "  -7.70%  " },
    @{ Row = 45; B = "Monero";           C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr";           D = "146.07"; E = "  -2.91%  " },
    @{ Row = 46; B = "Stellar";          C = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm";          D = "0.130";  E = "  +4.39%  " },
    @{ Row = 47; B = "OKB";              C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb";              D = "43.31";  E = "  -0.71%  " }
)

foreach ($u in $rowSwaps) {
    $r = $u.Row
    Set-TextCell $r 2 $u.B
    Set-TextCell $r 3 $u.C
    Set-TextCell $r 4 $u.D
    Set-TextCell $r 5 $u.E
}
